$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.254529666666667
$ws.Range("H2").Value = 3.763589
$ws.Range("I2").Value = 0.01270475613604028
$ws.Range("J2").Value = 0.01270475613604028
$ws.Range("M2").Value = 2.598166333333333
$ws.Range("N2").Value = 7.794499
$ws.Range("O2").Value = 0.3466013321552429
$ws.Range("P2").Value = 0.3466013321552429
$ws.Range("Q2").Value = 3.259476744101223
$ws.Range("R2").Value = 29.335290696911
$ws.Range("S2").Value = 0.004403485401459056
$ws.Range("T2").Value = 0.004403485401459056
$ws.Range("G3").Value = 1.254529666666667
$ws.Range("H3").Value = 3.763589
$ws.Range("I3").Value = 0.01270475613604028
$ws.Range("J3").Value = 0.01270475613604028
$ws.Range("M3").Value = 4.333403333333333
$ws.Range("O3").Value = 0.5780859172985858
$ws.Range("P3").Value = 0.5780859172985858
$ws.Range("Q3").Value = 5.436383039298889
$ws.Range("R3").Value = 48.92744735369
$ws.Range("S3").Value = 0.007344440604957681
$ws.Range("T3").Value = 0.007344440604957679
$ws.Range("G4").Value = 1.254529666666667
$ws.Range("H4").Value = 3.763589
$ws.Range("I4").Value = 0.01270475613604028
$ws.Range("J4").Value = 0.01270475613604028
$ws.Range("M4").Value = 0.4692043333333333
$ws.Range("N4").Value = 1.407613
$ws.Range("O4").Value = 0.06259293136852516
$ws.Range("P4").Value = 0.06259293136852516
$ws.Range("Q4").Value = 0.5886307558952223
$ws.Range("R4").Value = 5.297676803057001
$ws.Range("S4").Value = 0.000795227928877018
$ws.Range("T4").Value = 0.0007952279288770179
$ws.Range("G5").Value = 1.254529666666667
$ws.Range("H5").Value = 3.763589
$ws.Range("I5").Value = 0.01270475613604028
$ws.Range("J5").Value = 0.01270475613604028
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.09534933333333333
$ws.Range("N5").Value = 0.286048
$ws.Range("O5").Value = 0.01271981917764605
$ws.Range("P5").Value = 0.01271981917764604
$ws.Range("Q5").Value = 0.1196185673635556
$ws.Range("R5").Value = 1.076567106272
$ws.Range("S5").Value = 0.0001616022007465214
$ws.Range("T5").Value = 0.0001616022007465214
$ws.Range("I6").Value = 0.9734152842234517
$ws.Range("J6").Value = 0.9734152842234516
$ws.Range("M6").Value = 2.598166333333333
$ws.Range("N6").Value = 7.794499
$ws.Range("O6").Value = 0.3466013321552429
$ws.Range("P6").Value = 0.3466013321552429
$ws.Range("Q6").Value = 249.7351737652403
$ws.Range("R6").Value = 2247.616563887163
$ws.Range("S6").Value = 0.3373870342521227
$ws.Range("T6").Value = 0.3373870342521227
$ws.Range("I7").Value = 0.9734152842234517
$ws.Range("J7").Value = 0.9734152842234516
$ws.Range("M7").Value = 4.333403333333333
$ws.Range("O7").Value = 0.5780859172985858
$ws.Range("P7").Value = 0.5780859172985858
$ws.Range("Q7").Value = 416.5257707178633
$ws.Range("R7").Value = 3748.73193646077
$ws.Range("S7").Value = 0.5627176674927777
$ws.Range("T7").Value = 0.5627176674927776
$ws.Range("I8").Value = 0.9734152842234517
$ws.Range("J8").Value = 0.9734152842234516
$ws.Range("M8").Value = 0.4692043333333333
$ws.Range("N8").Value = 1.407613
$ws.Range("O8").Value = 0.06259293136852516
$ws.Range("P8").Value = 0.06259293136852516
$ws.Range("Q8").Value = 45.09981682584234
$ws.Range("R8").Value = 405.8983514325811
$ws.Range("S8").Value = 0.06092891607847192
$ws.Range("T8").Value = 0.06092891607847192
$ws.Range("I9").Value = 0.9734152842234517
$ws.Range("J9").Value = 0.9734152842234516
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.09534933333333333
$ws.Range("N9").Value = 0.286048
$ws.Range("O9").Value = 0.01271981917764605
$ws.Range("P9").Value = 0.01271981917764604
$ws.Range("Q9").Value = 9.164956847797333
$ws.Range("R9").Value = 82.484611630176
$ws.Range("S9").Value = 0.01238166640007924
$ws.Range("T9").Value = 0.01238166640007923
$ws.Range("G10").Value = 1.151276666666667
$ws.Range("H10").Value = 3.45383
$ws.Range("I10").Value = 0.01165910195968263
$ws.Range("J10").Value = 0.01165910195968263
$ws.Range("M10").Value = 2.598166333333333
$ws.Range("N10").Value = 7.794499
$ws.Range("O10").Value = 0.3466013321552429
$ws.Range("P10").Value = 0.3466013321552429
$ws.Range("Q10").Value = 2.991208275685556
$ws.Range("R10").Value = 26.92087448117
$ws.Range("S10").Value = 0.004041060270959803
$ws.Range("T10").Value = 0.004041060270959802
$ws.Range("G11").Value = 1.151276666666667
$ws.Range("H11").Value = 3.45383
$ws.Range("I11").Value = 0.01165910195968263
$ws.Range("J11").Value = 0.01165910195968263
$ws.Range("M11").Value = 4.333403333333333
$ws.Range("O11").Value = 0.5780859172985858
$ws.Range("P11").Value = 0.5780859172985858
$ws.Range("Q11").Value = 4.988946144922222
$ws.Range("R11").Value = 44.90051530429999
$ws.Range("S11").Value = 0.006739962651240873
$ws.Range("T11").Value = 0.006739962651240872
$ws.Range("G12").Value = 1.151276666666667
$ws.Range("H12").Value = 3.45383
$ws.Range("I12").Value = 0.01165910195968263
$ws.Range("J12").Value = 0.01165910195968263
$ws.Range("M12").Value = 0.4692043333333333
$ws.Range("N12").Value = 1.407613
$ws.Range("O12").Value = 0.06259293136852516
$ws.Range("P12").Value = 0.06259293136852516
$ws.Range("Q12").Value = 0.5401840008655556
$ws.Range("R12").Value = 4.86165600779
$ws.Range("S12").Value = 0.0007297773687810521
$ws.Range("T12").Value = 0.0007297773687810519
$ws.Range("G13").Value = 1.151276666666667
$ws.Range("H13").Value = 3.45383
$ws.Range("I13").Value = 0.01165910195968263
$ws.Range("J13").Value = 0.01165910195968263
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.09534933333333333
$ws.Range("N13").Value = 0.286048
$ws.Range("O13").Value = 0.01271981917764605
$ws.Range("P13").Value = 0.01271981917764604
$ws.Range("Q13").Value = 0.1097734626488889
$ws.Range("R13").Value = 0.9879611638399999
$ws.Range("S13").Value = 0.0001483016687009017
$ws.Range("T13").Value = 0.0001483016687009017
$ws.Range("G14").Value = 0.2192983333333333
$ws.Range("H14").Value = 0.6578949999999999
$ws.Range("I14").Value = 0.002220857680825461
$ws.Range("J14").Value = 0.002220857680825461
$ws.Range("M14").Value = 2.598166333333333
$ws.Range("N14").Value = 7.794499
$ws.Range("O14").Value = 0.3466013321552429
$ws.Range("P14").Value = 0.3466013321552429
$ws.Range("Q14").Value = 0.5697735466227777
$ws.Range("R14").Value = 5.127961919604999
$ws.Range("S14").Value = 0.0007697522307013082
$ws.Range("T14").Value = 0.000769752230701308
$ws.Range("G15").Value = 0.2192983333333333
$ws.Range("H15").Value = 0.6578949999999999
$ws.Range("I15").Value = 0.002220857680825461
$ws.Range("J15").Value = 0.002220857680825461
$ws.Range("M15").Value = 4.333403333333333
$ws.Range("O15").Value = 0.5780859172985858
$ws.Range("P15").Value = 0.5780859172985858
$ws.Range("Q15").Value = 0.9503081286611108
$ws.Range("R15").Value = 8.552773157949998
$ws.Range("S15").Value = 0.001283846549609597
$ws.Range("T15").Value = 0.001283846549609596
$ws.Range("G16").Value = 0.2192983333333333
$ws.Range("H16").Value = 0.6578949999999999
$ws.Range("I16").Value = 0.002220857680825461
$ws.Range("J16").Value = 0.002220857680825461
$ws.Range("M16").Value = 0.4692043333333333
$ws.Range("N16").Value = 1.407613
$ws.Range("O16").Value = 0.06259293136852516
$ws.Range("P16").Value = 0.06259293136852516
$ws.Range("Q16").Value = 0.1028957282927778
$ws.Range("R16").Value = 0.9260615546349998
$ws.Range("S16").Value = 0.0001390099923951701
$ws.Range("T16").Value = 0.00013900999239517
$ws.Range("G17").Value = 0.2192983333333333
$ws.Range("H17").Value = 0.6578949999999999
$ws.Range("I17").Value = 0.002220857680825461
$ws.Range("J17").Value = 0.002220857680825461
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.09534933333333333
$ws.Range("N17").Value = 0.286048
$ws.Range("O17").Value = 0.01271981917764605
$ws.Range("P17").Value = 0.01271981917764604
$ws.Range("Q17").Value = 0.02090994988444444
$ws.Range("R17").Value = 0.18818954896
$ws.Range("S17").Value = 0.000028248908119386219263326904
$ws.Range("T17").Value = 0.000028248908119386209098931537
